# Update cryptocurrency price/volume data per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.873.86'
$ws.Range('D3').Value = '2.044.68'
$ws.Range('E3').Value = '  -0.38%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').Value = '''245.20'
$ws.Range('E5').Value = '  -1.76%  '
$ws.Range('D6').Value = '''0.653'
$ws.Range('E6').Value = '  -2.28%  '
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('B8').Value = 'Solana'
$ws.Range('C8').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D8').Value = '''57.37'
$ws.Range('E8').Value = '  -4.36%  '
$ws.Range('D9').Value = '''0.368'
$ws.Range('E9').Value = '  -5.36%  '
$ws.Range('D10').Value = '''0.0772'
$ws.Range('E10').Value = '  -2.79%  '
$ws.Range('E11').Value = '  +1.41%  '
$ws.Range('D12').Value = '''15.14'
$ws.Range('E12').Value = '  -5.99%  '
$ws.Range('D13').Value = '''0.867'
$ws.Range('E13').Value = '  +3.69%  '
$ws.Range('D14').Value = '2.346.65'
$ws.Range('E14').Value = '  -0.21%  '
$ws.Range('E15').Value = '  -4.05%  '
$ws.Range('D16').Value = '2.011.29'
$ws.Range('E16').Value = '  -2.05%  '
$ws.Range('D17').Value = '''17.91'
$ws.Range('E17').Value = '  -2.55%  '
$ws.Range('D18').Value = '36.803.37'
$ws.Range('E18').Value = '  -0.73%  '
$ws.Range('D19').Value = '''73.21'
$ws.Range('E19').Value = '  -3.61%  '
$ws.Range('D20').Value = '0.0₃0882'
$ws.Range('E20').Value = '  -2.57%  '
$ws.Range('D21').Value = '''5.36'
$ws.Range('E21').Value = '  -0.67%  '
$ws.Range('D22').Value = '''235.36'
$ws.Range('E22').Value = '  -1.27%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('D25').Value = '''10.19'
$ws.Range('E25').Value = '  +7.95%  '
$ws.Range('E26').Value = '  -1.42%  '
$ws.Range('D27').Value = '''168.34'
$ws.Range('E27').Value = '  -0.77%  '
$ws.Range('D28').Value = '''19.87'
$ws.Range('E28').Value = '  -1.77%  '
$ws.Range('E29').Value = '  +13.03%  '
$ws.Range('E30').Value = '  -2.56%  '
$ws.Range('E31').Value = '  -4.39%  '
$ws.Range('E32').Value = '  +2.25%  '
$ws.Range('E33').Value = '  -4.16%  '
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('D35').Value = '''2.31'
$ws.Range('E35').Value = '  +3.09%  '
$ws.Range('E36').Value = '  +4.44%  '
$ws.Range('D37').Value = '''0.0821'
$ws.Range('E37').Value = '  -7.95%  '
$ws.Range('E38').Value = '  -3.09%  '
$ws.Range('E39').Value = '  -3.13%  '
$ws.Range('E40').Value = '  -5.41%  '
$ws.Range('E41').Value = '  -1.18%  '
$ws.Range('D42').Value = '''1.13'
$ws.Range('E42').Value = '  -0.69%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Value = '''96.18'
$ws.Range('E43').Value = '  -1.27%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').Value = '''0.0935'
$ws.Range('E44').Value = '  -14.42%  '
$ws.Range('D45').Value = '''16.78'
$ws.Range('E45').Value = '  -4.99%  '
$ws.Range('D46').Value = '1.297.39'
$ws.Range('E46').Value = '  +0.17%  '
$ws.Range('D47').Value = '''2.34'
$ws.Range('E47').Value = '  -7.22%  '
$ws.Range('E48').Value = '  -1.19%  '
$ws.Range('D49').Value = '''6.71'
$ws.Range('E49').Value = '  -2.46%  '
$ws.Range('D50').Value = '2.230.65'
$ws.Range('E50').Value = '  -0.40%  '
$ws.Range('D51').Value = '''44.25'
$ws.Range('E51').Value = '  -0.27%  '
